$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new headers in I1:J1, copying H1's formatting (bold, centered,
# bordered header style) so the new header cells match the existing ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Row 2 is a special case: I0 = 4, IF = 4 (does not mirror the IP column).
$ws.Cells.Item(2, 9).Value2 = 4
$ws.Cells.Item(2, 10).Value2 = 4

# Remaining data rows (3-28): I0 is always 1, and IF mirrors the IP (H) value.
for ($r = 3; $r -le 28; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $ipValue
}
